$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text format (not auto-converted to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.703.14"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").Value = "1.770.91"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "304.55"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("D7").Value = "0.4373"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("D8").Value = "0.3635"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").Value = "0.07174"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "0.8360"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").Value = "20.18"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("D12").Value = "1.761.36"
$ws.Range("E12").Value = "  -10.89%  "
$ws.Range("D13").Value = "5.241"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").Value = "6.340"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "0.06804"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "79.30"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "0.000008691"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -3.09%  "
$ws.Range("D21").Value = "26.477.29"
$ws.Range("E21").Value = "  -5.34%  "
$ws.Range("D22").Value = "5.012"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("D24").Value = "1.942.50"
$ws.Range("E24").Value = "  -10.53%  "
$ws.Range("D25").Value = "1.901"
$ws.Range("E25").Value = "  -4.68%  "
$ws.Range("D26").Value = "153.22"
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").Value = "18.17"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("D28").Value = "5.052"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").Value = "114.32"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "1.637"
$ws.Range("E30").Value = "  -10.35%  "
$ws.Range("D31").Value = "0.08981"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").Value = "0.7169"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("D33").Value = "4.317"
$ws.Range("E33").Value = "  -4.80%  "
$ws.Range("D34").Value = "2.798"
$ws.Range("E34").Value = "  -6.28%  "
$ws.Range("D35").Value = "1.082"
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "1.074"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").Value = "0.05097"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").Value = "0.01887"
$ws.Range("E39").Value = "  -2.37%  "
$ws.Range("D40").Value = "0.4915"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").Value = "0.1605"
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("D42").Value = "2.539"
$ws.Range("E42").Value = "  -9.25%  "
$ws.Range("D43").Value = "6.140"
$ws.Range("E43").Value = "  -5.97%  "
$ws.Range("D44").Value = "7.900"
$ws.Range("E44").Value = "  -4.85%  "
$ws.Range("D45").Value = "104.71"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "10.06"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").Value = "0.06212"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").Value = "0.4482"
$ws.Range("E49").Value = "  -4.10%  "
$ws.Range("D50").Value = "1.574"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").Value = "1.701"
$ws.Range("E51").Value = "  -0.66%  "
